$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D14: mark result "failed"
$ws.Range("D14").Value = "failed"

# D21: mark result "cleared"
$ws.Range("D21").Value = "cleared"

# New row 22 - copy formatting from row 21 first, then overwrite values
$ws.Range("A21:C21").Copy($ws.Range("A22:C22"))
$ws.Range("A22").Value = 45930
$ws.Range("B22").Value = "Sharekhan , Kanjur"
$ws.Range("C22").Value = "code for producer consumer problem and consumer will be process by two threads, print 1 to 10 using 5 threads, http vs websocket, tcp"
$ws.Rows.Item(22).RowHeight = 45

$ws.Range("C27").Select() | Out-Null
